# Översikt KUNGSÖR.xlsx update:
#  - Every data row's "Förändrad" date (column C) moves from 46063 to 46064.
#  - The data rows (5..79) are re-ordered according to a refreshed export
#    (same logical records, new row order) - captured below as oldRow -> newRow.
#
# Strategy: snapshot every data cell (A..Z) for rows 2..79 from the ORIGINAL
# sheet into memory (value or formula), then clear A5:Z79 and rewrite every
# row from the snapshot at its new location, bumping column C by 1 day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 79
$lastCol = 26   # column Z

# ---- 1. Snapshot all cells (A..Z) for rows 2..79 -----------------------
$snapshot = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $rowData = @{}
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.HasFormula) {
            $rowData[$c] = @{ "f" = $true; "v" = $cell.Formula }
        } else {
            $v = $cell.Value2
            if ($v -ne $null) {
                $rowData[$c] = @{ "f" = $false; "v" = $v }
            }
        }
    }
    $snapshot[$r] = $rowData
}

# ---- 2. New-row -> old-row mapping for the re-ordered block (rows 5-79) ---
$rowMap = @{
    5 = 6
    6 = 5
    7 = 7
    8 = 8
    9 = 9
    10 = 10
    11 = 11
    12 = 12
    13 = 13
    14 = 14
    15 = 15
    16 = 16
    17 = 17
    18 = 18
    19 = 19
    20 = 58
    21 = 59
    22 = 38
    23 = 29
    24 = 30
    25 = 31
    26 = 66
    27 = 67
    28 = 41
    29 = 73
    30 = 44
    31 = 68
    32 = 34
    33 = 63
    34 = 40
    35 = 64
    36 = 45
    37 = 74
    38 = 62
    39 = 76
    40 = 75
    41 = 78
    42 = 79
    43 = 77
    44 = 65
    45 = 32
    46 = 57
    47 = 70
    48 = 60
    49 = 69
    50 = 36
    51 = 37
    52 = 21
    53 = 35
    54 = 23
    55 = 24
    56 = 55
    57 = 56
    58 = 71
    59 = 22
    60 = 28
    61 = 39
    62 = 25
    63 = 52
    64 = 53
    65 = 54
    66 = 47
    67 = 20
    68 = 26
    69 = 72
    70 = 50
    71 = 27
    72 = 42
    73 = 49
    74 = 61
    75 = 51
    76 = 46
    77 = 33
    78 = 48
    79 = 43
}

# ---- 3. Clear the block that gets reshuffled (rows 5..79, cols A..Z) ----
$ws.Range("A5:Z79").ClearContents()

# ---- 4. Rewrite rows 2..4 in place, bumping column C only --------------
for ($r = $firstDataRow; $r -le 4; $r++) {
    $cCell = $snapshot[$r][3]
    if ($cCell -ne $null) {
        $ws.Cells.Item($r, 3).Value = [double]$cCell["v"] + 1
    }
}

# ---- 5. Rewrite rows 5..79 from the snapshot using the new row mapping --
foreach ($newRow in $rowMap.Keys) {
    $oldRow = $rowMap[$newRow]
    $rowData = $snapshot[$oldRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $entry = $rowData[$c]
        if ($entry -eq $null) {
            continue
        }
        $destCell = $ws.Cells.Item($newRow, $c)
        if ($entry["f"]) {
            $destCell.Formula = $entry["v"]
        } else {
            if ($c -eq 3) {
                $destCell.Value = [double]$entry["v"] + 1
            } else {
                $destCell.Value = $entry["v"]
            }
        }
    }
}
